{"js": "// \"New dream, new song\" \u2014 the red annotation's branch name changes from\n// \"alternate\" to \"main\": (This is a change \u2013 Version for branch alternate)\n// becomes (This is a change \u2013 Version for branch main), with \"main\" split\n// into its own run (mirroring a type-over-selection edit in Word).\n\n// Locate the whole parenthetical run by its current (pre-edit) text so the\n// edit is anchored precisely, regardless of surrounding content.\nconst target = \"(This is a change \\u2013 Version for branch alternate)\";\nconst results = context.document.body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text to edit.\");\n}\n\nconst targetRange = results.items[0];\n\n// Replace the matched range with three runs that keep the original red\n// color but isolate the new word \"main\" in its own run:\n//   \"(This is a change \u2013 Version for branch \" | \"main\" | \")\"\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:rPr><w:color w:val=\"C00000\"/></w:rPr><w:t xml:space=\"preserve\">(This is a change \\u2013 Version for branch </w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"C00000\"/></w:rPr><w:t>main</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"C00000\"/></w:rPr><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"New dream, new song\" - the red annotation's branch name changes from\n# \"alternate\" to \"main\": (This is a change \u2013 Version for branch alternate)\n# becomes (This is a change \u2013 Version for branch main), with \"main\" split\n# into its own run (mirroring a type-over-selection edit in Word).\n\n$d = $word.ActiveDocument\n\n# Locate the whole parenthetical run by its current (pre-edit) text so the\n# edit is anchored precisely, regardless of surrounding content.\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"(This is a change \u2013 Version for branch alternate)\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find target text to edit.\"\n}\n\n# Replace the matched range with three runs that keep the original red\n# color but isolate the new word \"main\" in its own run:\n#   \"(This is a change \u2013 Version for branch \" | \"main\" | \")\"\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:rPr><w:color w:val=\"C00000\"/></w:rPr><w:t xml:space=\"preserve\">(This is a change \u2013 Version for branch </w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"C00000\"/></w:rPr><w:t>main</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"C00000\"/></w:rPr><w:t>)</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$range.InsertXML($ooxml)\n"}
